$wb = $excel.ActiveWorkbook

# Report generated for archive: the former "Ready for handoff" status is now
# "In Translation" everywhere it appears (Overview zh-cn/de-de columns plus
# the per-language "Status" columns).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value()
        if ("Ready for handoff" -eq $val) {
            $cell.Value = "In Translation"
        }
    }
}

# The shorter status text means the Status columns can be narrower.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
